$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.918.05'
$ws.Range('E2').Value = '  -1.65%  '
$ws.Range('D3').Value = '2.450.69'
$ws.Range('E3').Value = '  -2.79%  '
$ws.Range('D4').Value = "'0.999"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').Value = "'578.27"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.02%  '
$ws.Range('D6').Value = "'165.44"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -6.07%  '
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('D8').Value = "'0.510"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -3.61%  '
$ws.Range('D9').Value = '2.451.16'
$ws.Range('E9').Value = '  -2.70%  '
$ws.Range('D10').Value = "'0.133"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -4.69%  '
$ws.Range('D11').Value = "'0.164"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.86%  '
$ws.Range('D12').Value = "'0.332"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -3.52%  '
$ws.Range('D13').Value = "'4.85"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -4.96%  '
$ws.Range('B14').Value = 'Avalanche'
$ws.Range('C14').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D14').Value = "'25.18"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -5.22%  '
$ws.Range('B15').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C15').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D15').Value = '2.874.80'
$ws.Range('E15').Value = '  -3.45%  '
$ws.Range('D16').Value = '66.543.81'
$ws.Range('E16').Value = '  -2.19%  '
$ws.Range('D17').Value = "'0.0000167"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -6.20%  '
$ws.Range('D18').Value = '2.435.89'
$ws.Range('E18').Value = '  -3.82%  '
$ws.Range('D19').Value = "'7.70"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -4.51%  '
$ws.Range('D20').Value = "'353.06"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.07%  '
$ws.Range('D21').Value = "'9.85"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -17.51%  '
$ws.Range('D22').Value = "'4.04"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.14%  '
$ws.Range('E23').Value = '  -0.13%  '
$ws.Range('D24').Value = "'68.94"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -3.12%  '
$ws.Range('D25').Value = "'4.23"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -9.18%  '
$ws.Range('D26').Value = "'1.74"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -9.74%  '
$ws.Range('D27').Value = "'8.89"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -11.94%  '
$ws.Range('E28').Value = '  -0.20%  '
$ws.Range('D29').Value = '2.558.90'
$ws.Range('E29').Value = '  -3.20%  '
$ws.Range('D30').Value = '0.0₃0898'
$ws.Range('E30').Value = '  -8.61%  '
$ws.Range('D31').Value = "'506.84"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -4.58%  '
$ws.Range('D32').Value = "'7.82"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -5.91%  '
$ws.Range('D33').Value = "'1.77"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -7.06%  '
$ws.Range('D34').Value = "'1.22"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -8.18%  '
$ws.Range('D35').Value = "'0.999"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.06%  '
$ws.Range('E36').Value = '  +0.69%  '
$ws.Range('D37').Value = "'0.116"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -10.16%  '
$ws.Range('D38').Value = "'18.56"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.61%  '
$ws.Range('D39').Value = "'18.36"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.01%  '
$ws.Range('D40').Value = "'1.34"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -7.47%  '
$ws.Range('E41').Value = '  +0.06%  '
$ws.Range('D42').Value = "'1.66"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -7.62%  '
$ws.Range('D43').Value = "'0.325"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -7.40%  '
$ws.Range('D44').Value = "'4.73"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -8.08%  '
$ws.Range('D45').Value = "'38.69"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.00%  '
$ws.Range('D46').Value = "'2.30"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -8.16%  '
$ws.Range('D47').Value = "'140.86"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -4.37%  '
$ws.Range('D48').Value = "'3.47"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -6.33%  '
$ws.Range('D49').Value = "'0.512"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -7.64%  '
$ws.Range('D50').Value = "'1.59"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -7.59%  '
$ws.Range('D51').Value = "'0.0731"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.76%  '
